# changes4polls.xlsx — fix the broken "poll_template_votebody" HTML template.
#
# The <input> tag's `name` attribute was emitted with a stray/doubled quote:
#     name=""poll_%POLL_ID%"
# which is invalid HTML. This corrects it to:
#     name="poll_%POLL_ID%"
#
# That template string lives in cell C3 (row "poll_template_votebody").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedVoteBody = '<li class="clearfix"> <label for="poll-answer-%POLL_ANSWER_ID%">%POLL_ANSWER%</label> <span><input type="%POLL_CHECKBOX_RADIO%" name="poll_%POLL_ID%" value="%POLL_ANSWER_ID%" id="poll-answer-%POLL_ANSWER_ID%"></span> </li>'

$ws.Range("C3").Value = $fixedVoteBody

# Leave the selection where the editor's cursor ended up after making the edit.
$ws.Range("D3").Select()
